$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 320-321 (existing rows 320-327 shift down to 322-329)
$ws.Rows("320:321").Insert()

# --- New row 320 ---
$ws.Range("A320").Value2 = 7
$ws.Range("B320").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C320").Value2 = "Ñuble"
$ws.Range("D320").Value2 = 44448
$ws.Range("E320").Value2 = 16
$ws.Range("F320").Value2 = 100112004
$ws.Range("G320").Value2 = "Cebolla"
$ws.Range("H320").Value2 = "Sin especificar"
$ws.Range("I320").Value2 = "1a (guarda)"
$ws.Range("J320").Value2 = 300
$ws.Range("K320").Value2 = 5500
$ws.Range("L320").Value2 = 6000
$ws.Range("M320").Value2 = 5750
$ws.Range("N320").Value2 = "$/malla 25 kilos"
$ws.Range("O320").Value2 = "Región del Maule"
$ws.Range("P320").Value2 = 230
$ws.Range("Q320").Value2 = 25
$ws.Range("R320").Value2 = "Hortaliza"

# --- New row 321 ---
$ws.Range("A321").Value2 = 7
$ws.Range("B321").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C321").Value2 = "Ñuble"
$ws.Range("D321").Value2 = 44448
$ws.Range("E321").Value2 = 16
$ws.Range("F321").Value2 = 100112004
$ws.Range("G321").Value2 = "Cebolla"
$ws.Range("H321").Value2 = "Sin especificar"
$ws.Range("I321").Value2 = "2a (guarda)"
$ws.Range("J321").Value2 = 150
$ws.Range("K321").Value2 = 4500
$ws.Range("L321").Value2 = 4500
$ws.Range("M321").Value2 = 4500
$ws.Range("N321").Value2 = "$/malla 25 kilos"
$ws.Range("O321").Value2 = "Región del Maule"
$ws.Range("P321").Value2 = 180
$ws.Range("Q321").Value2 = 25
$ws.Range("R321").Value2 = "Hortaliza"

# Ensure the date cells keep the date number format used by the rest of column D
$ws.Range("D320").NumberFormat = $ws.Range("D322").NumberFormat
$ws.Range("D321").NumberFormat = $ws.Range("D322").NumberFormat
